# Fruta / hortaliza, semanal
# Insert a new weekly data row before the current row 69, shifting the
# existing rows 69-71 down to 70-72, and populate the new row with the
# latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 69; rows 69-71 shift down to 70-72.
$ws.Rows(69).Insert()

# Fill in the new row 69 with this week's data.
$ws.Cells.Item(69, 1).Value = 8
$ws.Cells.Item(69, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(69, 3).Value = "Coquimbo"
$ws.Cells.Item(69, 4).Value = 44610
$ws.Cells.Item(69, 5).Value = 4
$ws.Cells.Item(69, 6).Value = 100112030
$ws.Cells.Item(69, 7).Value = "Poroto granado"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 400
$ws.Cells.Item(69, 11).Value = 30000
$ws.Cells.Item(69, 12).Value = 31000
$ws.Cells.Item(69, 13).Value = 30500
$ws.Cells.Item(69, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(69, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(69, 16).Value = 1220
$ws.Cells.Item(69, 17).Value = 25
$ws.Cells.Item(69, 18).Value = "Hortaliza"
